$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new test-case row at row 9 (existing rows 9-12 shift down to 10-13) ---
$ws.Rows.Item(9).Insert()

# Copy the formatting (styles/borders/alignment) of the row right below (old row 9,
# now shifted to row 10) into the freshly inserted blank row 9, so the new row matches
# the look of the other data rows exactly.
$ws.Range("B10:H10").Copy()
$ws.Range("B9:H9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new test case's data (TC #6: invalid-choice validation)
$ws.Cells.Item(9, 2).Value = 6
$ws.Cells.Item(9, 4).Value = "Verify that the application does not accept any other choice apart from 1,2 & 3."
$ws.Cells.Item(9, 5).Value = "P1"
$ws.Cells.Item(9, 6).Value = "Apllication should not accept any other choice apart from 1,2 & 3 and should be showing an error message to the user."
$ws.Cells.Item(9, 7).Value = "Apllication does not accept any other choice apart from 1,2 & 3 and shows an error message as ""Please Enter a Valid Choice!""."
$ws.Cells.Item(9, 8).Value = "Pass"

# Row-height tweaks
$ws.Rows.Item(8).RowHeight = 36
$ws.Rows.Item(9).RowHeight = 54

# --- Column width tweaks ---
$ws.Columns.Item(4).ColumnWidth = 45
$ws.Columns.Item(7).ColumnWidth = 40.833333333333336

# --- View state: selection moves to O16 ---
$ws.Range("O16").Select()
